$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Stage the cell-format "templates" we need in a scratch area (column Z)
#    BEFORE we start overwriting B3:B20, so every paste-format below reuses
#    an existing fill/style instead of inventing duplicates.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats -> style "s6" (light yellow)

$ws.Range("B12").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # xlPasteFormats -> style "s4" (light blue)

$ws.Range("B17").Copy()
$ws.Range("Z3").PasteSpecial(-4122)   # xlPasteFormats -> style "s5" (green 92D050)

# Brand new fill colour used by the "Avatar p1" grouping (FF8FCC4F).
$ws.Range("Z4").Interior.Color = 0x4FCC8F

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Re-populate B3:B20 with the re-sorted feature list and re-apply the
#    matching fill for every row.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 3;  Text = "Audio control (users Pitch, volume)";               Style = "Z1" },
    @{ Row = 4;  Text = "Host room";                                          Style = "Z1" },
    @{ Row = 5;  Text = "Join room";                                          Style = "Z1" },
    @{ Row = 6;  Text = "Multiplayer speech";                                 Style = "Z1" },
    @{ Row = 7;  Text = "Hands/body to not phase through objects and wall?";  Style = "Z2" },
    @{ Row = 8;  Text = "Life like experience (First person walking motion)"; Style = "Z2" },
    @{ Row = 9;  Text = "Make a body";                                        Style = "Z2" },
    @{ Row = 10; Text = "Wrist menu";                                         Style = "Z4" },
    @{ Row = 11; Text = "Wrist Menu showing object to spawn";                 Style = "Z4" },
    @{ Row = 12; Text = "Adding ambient sound ";                              Style = "Z4" },
    @{ Row = 13; Text = "Make object (to be spawned into enviroment)";        Style = "Z4" },
    @{ Row = 14; Text = "Wrist Menu showing Audio control";                   Style = "Z4" },
    @{ Row = 15; Text = "Haptic feedback";                                    Style = "Z4" },
    @{ Row = 16; Text = "Make room";                                          Style = "Z4" },
    @{ Row = 17; Text = "Object Ownership";                                   Style = "Z4" },
    @{ Row = 18; Text = "Request oject client & server";                     Style = "Z4" },
    @{ Row = 19; Text = "make github";                                        Style = "Z3" },
    @{ Row = 20; Text = "Spawn objects into enviroment";                      Style = "Z3" }
)

foreach ($r in $rows) {
    $target = $ws.Cells.Item($r.Row, 2)   # column B
    $target.Value = $r.Text
    $ws.Range($r.Style).Copy()
    $target.PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Clear the scratch cells we used as format templates.
# ---------------------------------------------------------------------------
$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# 4. Misc view / layout tweaks captured in the diff.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 24
$ws.Application.ActiveWindow.Zoom = 156
$ws.Range("D18").Select()
